# Update yield and fertilizer calculation:
# - Add new unit label "[lbs P/ac]" to the shared strings table and use it
#   for the DAP fertilizer columns (L3/M3) on both the "inputs" and
#   "outputs" sheets (previously they incorrectly showed "[lbs N/ac]").
# - Remove the stray, empty formatted rows (17-37) that were left over on
#   the "outputs" sheet.
# - Leave "inputs" as the active/selected sheet with L3 selected, and
#   select M3 on "outputs".

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("inputs")
$wsOutputs = $wb.Worksheets.Item("outputs")

# Correct the unit label for the DAP (di-ammonium phosphate) fertilizer
# columns - these should be phosphorus pounds per acre, not nitrogen.
$wsInputs.Range("L3:M3").Value = "[lbs P/ac]"
$wsOutputs.Range("L3:M3").Value = "[lbs P/ac]"

# Drop the leftover blank formatted rows on the "outputs" sheet.
$wsOutputs.Rows("17:37").Delete()

# Restore the view/selection state.
$wsOutputs.Range("M3").Select() | Out-Null
$wsInputs.Activate() | Out-Null
$wsInputs.Range("L3").Select() | Out-Null
